$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 = "RCL Leader Election" entry - code edit confirmed for LCRLeader.rebeca
$ws.Range("D7").Value = "OK, Progress"

$ws.Range("E7").Value = "msgsrv -> constructor, removed comments"
$ws.Range("E7").Style = "Good"

$ws.Range("F7").Value = 28
$ws.Range("F7").Style = "Good"

$ws.Range("G7").Value = 49
$ws.Range("G7").Style = "Good"

# Update the selected cell shown in the sheet view
$ws.Range("D8").Select() | Out-Null
